$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold, border, centered) from the existing "sum" header (G1)
# onto the new "Save" header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add new numeric value for the "Save" column in the data row.
$ws.Range("H2").Value = 0
